$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 4-29 down to rows 6-31 (process bottom-up to avoid overwrite)
for ($r = 29; $r -ge 4; $r--) {
    $dstRow = $r + 2
    for ($c = 1; $c -le 23; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($dstRow, $c)
        $dstCell.Value = $srcCell.Value2
    }
}

# Fill in new row 4 (Holden) and row 5 (Rizzie Spiral) with fresh simulation data
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.8357975553510357
$ws.Range("D4").Value = 1.006826651687704
$ws.Range("E4").Value = 0.8911600224436456
$ws.Range("F4").Value = 1.52899200910204
$ws.Range("G4").Value = 0.8357975553510357
$ws.Range("H4").Value = 0.8830645045437154
$ws.Range("I4").Value = 1.16675846900188
$ws.Range("J4").Value = 0.8911600224436456
$ws.Range("K4").Value = 0.8911600224436456
$ws.Range("L4").Value = 0.855961491981313
$ws.Range("M4").Value = 0.9033784577399147
$ws.Range("N4").Value = 0.8911600224436456
$ws.Range("O4").Value = 1.52899200910204
$ws.Range("P4").Value = 1.182394782226538
$ws.Range("Q4").Value = 1.206028256822878
$ws.Range("R4").Value = 1.085316528965574
$ws.Range("S4").Value = 1.08261802299893
$ws.Range("T4").Value = 1.085316528965574
$ws.Range("U4").Value = 1.034753522860109
$ws.Range("V4").Value = 1.006034822776817
$ws.Range("W4").Value = 1.008992395231406

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 0.001385562194050447
$ws.Range("D5").Value = 0.5677721269685765
$ws.Range("E5").Value = 0.3138923632986431
$ws.Range("F5").Value = 2.485380016487212
$ws.Range("G5").Value = 0.001385562194050447
$ws.Range("H5").Value = 0.7208521572890533
$ws.Range("I5").Value = 1.788186260714938
$ws.Range("J5").Value = 0.3138923632986431
$ws.Range("K5").Value = 0.3138923632986431
$ws.Range("L5").Value = 0.1056032685047131
$ws.Range("M5").Value = 1.445914630321507
$ws.Range("N5").Value = 0.3138923632986431
$ws.Range("O5").Value = 2.485380016487212
$ws.Range("P5").Value = 1.243382789340631
$ws.Range("Q5").Value = 1.603116086888133
$ws.Range("R5").Value = 0.9335526473266352
$ws.Range("S5").Value = 1.069205911990105
$ws.Range("T5").Value = 0.9335526473266352
$ws.Range("U5").Value = 0.8803775248172397
$ws.Range("V5").Value = 0.7670804925135204
$ws.Range("W5").Value = 0.9286232982223367

# Fix column-A style for the two brand-new rows at the bottom (30, 31) to match existing look
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)

# Rename category label "Thomas Hex" -> "Matthies Hex" (now located at row 11 after the shift)
$ws.Range("B11").Value = "Matthies Hex"